# BWP Test Cases Added
# Updates the Katalon execution timestamps recorded in column B ("Date")
# of the test-result worksheets with the latest run's timestamps.

$wb = $excel.ActiveWorkbook

$updates = @{
    "PayNowCC" = @{
        "B2" = "Sat Nov 15 20:34:58 EST 2025"
        "B3" = "Sat Nov 15 20:35:35 EST 2025"
        "B4" = "Sat Nov 15 20:36:06 EST 2025"
        "B5" = "Sat Nov 15 20:36:35 EST 2025"
    }
    "PayNowCCSCF" = @{
        "B2" = "Sat Nov 15 20:37:05 EST 2025"
        "B3" = "Sat Nov 15 20:37:45 EST 2025"
        "B4" = "Sat Nov 15 20:38:25 EST 2025"
        "B5" = "Sat Nov 15 20:39:05 EST 2025"
    }
    "PayNowCCDCF" = @{
        "B2" = "Sat Nov 15 20:39:46 EST 2025"
        "B3" = "Sat Nov 15 20:40:25 EST 2025"
        "B4" = "Sat Nov 15 20:41:06 EST 2025"
        "B5" = "Sat Nov 15 20:41:47 EST 2025"
    }
    "NoModifyAmount" = @{
        "B2" = "Sat Nov 15 20:43:18 EST 2025"
    }
    "OverUnderPay" = @{
        "B2" = "Sat Nov 15 20:46:07 EST 2025"
        "B3" = "Sat Nov 15 20:46:26 EST 2025"
    }
    "NoOverPay" = @{
        "B2" = "Sat Nov 15 20:48:15 EST 2025"
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
